# Daily attendance processing - 2025-10-12 23:40:44
# Reverse the order of names/emails in the "Recorded By" column (G) for
# rows that have more than one recorder listed, skipping any row whose
# list includes "admin@admin.com".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $value = $cell.Value2

    if ($value -eq $null) { continue }
    if ($value -notmatch ",") { continue }
    if ($value -match "admin@admin.com") { continue }

    $parts = $value -split ",\s*"
    $reversed = $parts[($parts.Count - 1)..0]
    $cell.Value = [string]::Join(", ", $reversed)
}
